$d = $word.ActiveDocument
$t = $d.Tables(1)
$t.Cell(1,1).Range.Text = "81÷4=20, 1"
$t.Cell(1,2).Range.Text = "57÷4=14, 1"
$t.Cell(1,3).Range.Text = "55÷6=9, 1"
$t.Cell(1,4).Range.Text = "75÷5=15, 0"
$t.Cell(1,5).Range.Text = "74÷6=12, 2"
$t.Cell(5,1).Range.Text = "48÷2=24, 0"
$t.Cell(5,2).Range.Text = "54÷3=18, 0"
$t.Cell(5,3).Range.Text = "34÷9=3, 7"
$t.Cell(5,4).Range.Text = "85÷7=12, 1"
$t.Cell(5,5).Range.Text = "24÷3=8, 0"
$t.Cell(9,1).Range.Text = "34÷2=17, 0"
$t.Cell(9,2).Range.Text = "96÷3=32, 0"
$t.Cell(9,3).Range.Text = "17÷6=2, 5"
$t.Cell(9,4).Range.Text = "70÷3=23, 1"
$t.Cell(9,5).Range.Text = "36÷7=5, 1"
$t.Cell(13,1).Range.Text = "88÷5=17, 3"
$t.Cell(13,2).Range.Text = "48÷8=6, 0"
$t.Cell(13,3).Range.Text = "13÷5=2, 3"
$t.Cell(13,4).Range.Text = "92÷6=15, 2"
$t.Cell(13,5).Range.Text = "88÷5=17, 3"
$t.Cell(17,1).Range.Text = "85÷5=17, 0"
$t.Cell(17,2).Range.Text = "58÷7=8, 2"
$t.Cell(17,3).Range.Text = "82÷4=20, 2"
$t.Cell(17,4).Range.Text = "88÷3=29, 1"
$t.Cell(17,5).Range.Text = "56÷5=11, 1"
